$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1785.1
$ws.Range("I19").Value = 896.3333
$ws.Range("J19").Value = 2166
$ws.Range("K19").Value = 896.3333
$ws.Range("L19").Value = 2166
$ws.Range("M19").Value = -721.3333
$ws.Range("N19").Value = -2516
$ws.Range("H38").Value = 450.07144
$ws.Range("I38").Value = 40.2
$ws.Range("K38").Value = 120.6
$ws.Range("M38").Value = 251.4
$ws.Range("H43").Value = 4352.6
$ws.Range("I43").Value = 1600.5
$ws.Range("K43").Value = 1600.5
$ws.Range("M43").Value = -1531.5
$ws.Range("H137").Value = 41810.777
$ws.Range("I137").Value = 59443.805
$ws.Range("K137").Value = 178331.415
$ws.Range("M137").Value = -175781.415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 159.90909
$ws.Range("I5").Value = 63.25
$ws.Range("J5").Value = 215.14285
$ws.Range("K5").Value = 63.25
$ws.Range("L5").Value = 215.14285
$ws.Range("M5").Value = 48.75
$ws.Range("N5").Value = -439.14285
$ws.Range("H32").Value = 8042.5522
$ws.Range("I32").Value = 4554.041
$ws.Range("J32").Value = 19114.783
$ws.Range("K32").Value = 4554.041
$ws.Range("L32").Value = 19114.783
$ws.Range("M32").Value = -4267.041
$ws.Range("N32").Value = -19688.783
$ws.Range("H61").Value = 6246.2
$ws.Range("I61").Value = 6870.25
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 6870.25
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -6658.25
$ws.Range("N61").Value = -4174
$ws.Range("H74").Value = 51614.082
$ws.Range("I74").Value = 37252.24
$ws.Range("K74").Value = 37252.24
$ws.Range("M74").Value = -36378.24
$ws.Range("H77").Value = 51614.082
$ws.Range("I77").Value = 37252.24
$ws.Range("K77").Value = 186261.2
$ws.Range("M77").Value = -181893.2
$ws.Range("H97").Value = 1294802.8
$ws.Range("I97").Value = 1903148.8
$ws.Range("J97").Value = 2067.625
$ws.Range("K97").Value = 1903148.8
$ws.Range("L97").Value = 2067.625
$ws.Range("M97").Value = -1902652.8
$ws.Range("N97").Value = -3059.625
$ws.Range("H122").Value = 6947431.5
$ws.Range("I122").Value = 3995
$ws.Range("K122").Value = 11985
$ws.Range("M122").Value = -9535
$ws.Range("H132").Value = 4018.3333
$ws.Range("I132").Value = 2608
$ws.Range("J132").Value = 6133.8335
$ws.Range("K132").Value = 7824
$ws.Range("L132").Value = 18401.5005
$ws.Range("M132").Value = -5294
$ws.Range("N132").Value = -23461.5005
$ws.Range("H136").Value = 6246.2
$ws.Range("I136").Value = 6870.25
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 20610.75
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -18060.75
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 159.90909
$ws.Range("I4").Value = 63.25
$ws.Range("J4").Value = 215.14285
$ws.Range("K4").Value = 63.25
$ws.Range("L4").Value = 215.14285
$ws.Range("M4").Value = 51.75
$ws.Range("N4").Value = -445.14285
$ws.Range("H59").Value = 122996
$ws.Range("J59").Value = 122996
$ws.Range("L59").Value = 122996
$ws.Range("N59").Value = -124690

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8685.823
$ws.Range("I58").Value = 10726.272
$ws.Range("J58").Value = 4945
$ws.Range("K58").Value = 10726.272
$ws.Range("L58").Value = 4945
$ws.Range("M58").Value = -10523.272
$ws.Range("N58").Value = -5351
$ws.Range("H99").Value = 5665
$ws.Range("J99").Value = 5747.5
$ws.Range("L99").Value = 5747.5
$ws.Range("N99").Value = -8743.5
$ws.Range("H126").Value = 5665
$ws.Range("J126").Value = 5747.5
$ws.Range("L126").Value = 17242.5
$ws.Range("N126").Value = -22182.5
$ws.Range("H141").Value = 151882.53
$ws.Range("J141").Value = 151882.53
$ws.Range("L141").Value = 151882.53
$ws.Range("N141").Value = -162242.53

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 48432
$ws.Range("J5").Value = 167647.67
$ws.Range("L5").Value = 502943.01
$ws.Range("N5").Value = -503167.01
$ws.Range("H135").Value = 48432
$ws.Range("J135").Value = 167647.67
$ws.Range("L135").Value = 1508829.03
$ws.Range("N135").Value = -1513899.03

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H102").Value = 5037260.5
$ws.Range("I102").Value = 12347082
$ws.Range("J102").Value = 1382349.8
$ws.Range("K102").Value = 12347082
$ws.Range("L102").Value = 1382349.8
$ws.Range("M102").Value = -12345460
$ws.Range("N102").Value = -1385593.8
$ws.Range("H132").Value = 3517.3044
$ws.Range("I132").Value = 3452.2778
$ws.Range("J132").Value = 3751.4
$ws.Range("K132").Value = 10356.8334
$ws.Range("L132").Value = 11254.2
$ws.Range("M132").Value = -7826.8334
$ws.Range("N132").Value = -16314.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2177.182
$ws.Range("J16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("N16").Value = -5340
$ws.Range("H40").Value = 3610.9333
$ws.Range("I40").Value = 2670.652
$ws.Range("J40").Value = 6700.4287
$ws.Range("K40").Value = 2670.652
$ws.Range("L40").Value = 6700.4287
$ws.Range("M40").Value = -2534.652
$ws.Range("N40").Value = -6972.4287
$ws.Range("H46").Value = 4534.75
$ws.Range("I46").Value = 1083.3334
$ws.Range("K46").Value = 1083.3334
$ws.Range("M46").Value = -895.3334
$ws.Range("H54").Value = 43666.668
$ws.Range("J54").Value = 43666.668
$ws.Range("L54").Value = 43666.668
$ws.Range("N54").Value = -44954.668
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 3887.5557
$ws.Range("I122").Value = 2475.6155
$ws.Range("K122").Value = 7426.8465
$ws.Range("M122").Value = -4976.8465
$ws.Range("H133").Value = 119999.625
$ws.Range("J133").Value = 119999.625
$ws.Range("L133").Value = 119999.625
$ws.Range("N133").Value = -125059.625
